# Automatic map update (2025-11-13 08:03:47)
#
# The source data table on sheet "NEW" lost two records that are no longer
# relevant:
#   - Caso 6928 / "ALBARELLOS AV. 3101" (originally row 42)
#   - Caso -569 / "Palpa 2862"          (originally row 51)
#
# Removing those two rows shifts every following record up, which is why
# the sheet's used range shrinks from A1:R90 to A1:R88 and every row from
# 42 onward now shows the data that used to live two (or one, before the
# first deletion) rows further down.
#
# Delete from the bottom up so row numbers for the rows still to be
# removed don't shift out from under us.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(51).Delete()
$ws.Rows.Item(42).Delete()
